$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 19
$ws.Range("H19").Value = 871.25
$ws.Range("I19").Value = 644.25
$ws.Range("J19").Value = 1098.25
$ws.Range("K19").Value = 644.25
$ws.Range("L19").Value = 1098.25
$ws.Range("M19").Value = -469.25
$ws.Range("N19").Value = -1448.25

# row 40
$ws.Range("H40").Value = 1790.5
$ws.Range("I40").Value = 1731.8889
$ws.Range("J40").Value = 1966.3334
$ws.Range("K40").Value = 1731.8889
$ws.Range("L40").Value = 1966.3334
$ws.Range("M40").Value = -1556.8889
$ws.Range("N40").Value = -2316.3334

# row 70
$ws.Range("H70").Value = 92272.17999999999
$ws.Range("J70").Value = 92272.17999999999
$ws.Range("L70").Value = 276816.54
$ws.Range("N70").Value = -277356.54

# row 73
$ws.Range("H73").Value = 92272.17999999999
$ws.Range("J73").Value = 92272.17999999999
$ws.Range("L73").Value = 276816.54
$ws.Range("N73").Value = -278688.54

# row 94
$ws.Range("H94").Value = 991.75
$ws.Range("I94").Value = 991.75
$ws.Range("K94").Value = 991.75
$ws.Range("M94").Value = -540.75

# row 99
$ws.Range("H99").Value = 1734.4
$ws.Range("I99").Value = 2237.3333
$ws.Range("J99").Value = 980
$ws.Range("K99").Value = 6711.999899999999
$ws.Range("L99").Value = 2940
$ws.Range("M99").Value = -5213.999899999999
$ws.Range("N99").Value = -5936

# row 101
$ws.Range("H101").Value = 20000380
$ws.Range("I101").Value = 33333500
$ws.Range("J101").Value = 699.5
$ws.Range("K101").Value = 100000500
$ws.Range("L101").Value = 2098.5
$ws.Range("M101").Value = -99998878
$ws.Range("N101").Value = -5342.5

# row 112
$ws.Range("H112").Value = 2808.5
$ws.Range("J112").Value = 3128
$ws.Range("L112").Value = 9384
$ws.Range("N112").Value = -11600

# row 116
$ws.Range("H116").Value = 6500
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 6500
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 6500
$ws.Range("M116").ClearContents()
$ws.Range("N116").Value = -13384

# row 129
$ws.Range("H129").Value = 2316.0833
$ws.Range("I129").Value = 1259.6
$ws.Range("J129").Value = 3070.7144
$ws.Range("K129").Value = 3778.8
$ws.Range("L129").Value = 9212.143199999999
$ws.Range("M129").Value = 1221.2
$ws.Range("N129").Value = -19212.1432

# row 138
$ws.Range("H138").Value = 2705.0667
$ws.Range("J138").Value = 2922.6
$ws.Range("L138").Value = 8767.799999999999
$ws.Range("N138").Value = -19047.8

$ws = $wb.Worksheets.Item("ARM")
# row 21
$ws.Range("H21").Value = 5555
$ws.Range("J21").Value = 5555
$ws.Range("L21").Value = 5555
$ws.Range("N21").Value = -6303

# row 61
$ws.Range("H61").Value = 2044.72
$ws.Range("I61").Value = 2044.72
$ws.Range("K61").Value = 2044.72
$ws.Range("M61").Value = -1832.72

# row 88
$ws.Range("H88").Value = 604
$ws.Range("I88").Value = 491
$ws.Range("J88").Value = 668.5714
$ws.Range("K88").Value = 491
$ws.Range("L88").Value = 668.5714
$ws.Range("M88").Value = -85
$ws.Range("N88").Value = -1480.5714

# row 91
$ws.Range("H91").Value = 604
$ws.Range("I91").Value = 491
$ws.Range("J91").Value = 668.5714
$ws.Range("K91").Value = 491
$ws.Range("L91").Value = 668.5714
$ws.Range("M91").Value = 913
$ws.Range("N91").Value = -3476.5714

# row 122
$ws.Range("H122").Value = 2351.15
$ws.Range("I122").Value = 2157.9092
$ws.Range("K122").Value = 6473.7276
$ws.Range("M122").Value = -4023.7276

# row 125
$ws.Range("H125").Value = 71000
$ws.Range("J125").Value = 71000
$ws.Range("L125").Value = 71000
$ws.Range("N125").Value = -80840

# row 131
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()

# row 132
$ws.Range("H132").Value = 412
$ws.Range("I132").Value = 412
$ws.Range("K132").Value = 1236
$ws.Range("M132").Value = 1294

# row 136
$ws.Range("H136").Value = 2044.72
$ws.Range("I136").Value = 2044.72
$ws.Range("K136").Value = 6134.16
$ws.Range("M136").Value = -3584.16

# row 139
$ws.Range("H139").Value = 77454.664
$ws.Range("J139").Value = 75857
$ws.Range("L139").Value = 75857
$ws.Range("N139").Value = -86137

$ws = $wb.Worksheets.Item("BSM")
# row 12
$ws.Range("H12").Value = 495
$ws.Range("I12").Value = 495
$ws.Range("K12").Value = 495
$ws.Range("M12").Value = -327

# row 86
$ws.Range("H86").Value = 1494
$ws.Range("I86").Value = 1408.5
$ws.Range("K86").Value = 1408.5
$ws.Range("M86").Value = -285.5

# row 89
$ws.Range("H89").Value = 1494
$ws.Range("I89").Value = 1408.5
$ws.Range("K89").Value = 7042.5
$ws.Range("M89").Value = -1426.5

# row 94
$ws.Range("H94").Value = 664.13635
$ws.Range("I94").Value = 604.41174
$ws.Range("K94").Value = 604.41174
$ws.Range("M94").Value = -153.41174

# row 99
$ws.Range("H99").Value = 1294.7222
$ws.Range("I99").Value = 1093.5333
$ws.Range("K99").Value = 1093.5333
$ws.Range("M99").Value = 404.4666999999999

$ws = $wb.Worksheets.Item("CRP")
# row 31
$ws.Range("H31").Value = 2210.2
$ws.Range("I31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("M31").ClearContents()

# row 34
$ws.Range("H34").Value = 2210.2
$ws.Range("I34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("M34").ClearContents()

# row 99
$ws.Range("H99").Value = 1375
$ws.Range("I99").Value = 1375
$ws.Range("K99").Value = 1375
$ws.Range("M99").Value = 123

# row 105
$ws.Range("H105").Value = 3253.75
$ws.Range("I105").Value = 2804.8572
$ws.Range("J105").Value = 3882.2
$ws.Range("K105").Value = 2804.8572
$ws.Range("L105").Value = 3882.2
$ws.Range("M105").Value = -1057.8572
$ws.Range("N105").Value = -7376.2

# row 126
$ws.Range("H126").Value = 1375
$ws.Range("I126").Value = 1375
$ws.Range("K126").Value = 4125
$ws.Range("M126").Value = -1655

# row 140
$ws.Range("H140").Value = 79000
$ws.Range("J140").Value = 79000
$ws.Range("L140").Value = 79000
$ws.Range("N140").Value = -89360

# row 141
$ws.Range("H141").Value = 35923.066
$ws.Range("J141").Value = 35275.5
$ws.Range("L141").Value = 35275.5
$ws.Range("N141").Value = -45635.5

$ws = $wb.Worksheets.Item("CUL")
# row 69
$ws.Range("H69").Value = 2893.6
$ws.Range("I69").Value = 1991
$ws.Range("J69").Value = 2993.889
$ws.Range("K69").Value = 5973
$ws.Range("L69").Value = 8981.667000000001
$ws.Range("M69").Value = -5162
$ws.Range("N69").Value = -10603.667

# row 72
$ws.Range("H72").Value = 2893.6
$ws.Range("I72").Value = 1991
$ws.Range("J72").Value = 2993.889
$ws.Range("K72").Value = 17919
$ws.Range("L72").Value = 26945.001
$ws.Range("M72").Value = -13863
$ws.Range("N72").Value = -35057.001

# row 75
$ws.Range("H75").Value = 10061
$ws.Range("I75").Value = 1804
$ws.Range("K75").Value = 5412
$ws.Range("M75").Value = -4414

# row 78
$ws.Range("H78").Value = 10061
$ws.Range("I78").Value = 1804
$ws.Range("K78").Value = 16236
$ws.Range("M78").Value = -11244

# row 108
$ws.Range("H108").Value = 2209
$ws.Range("I108").Value = 2209
$ws.Range("K108").Value = 6627
$ws.Range("M108").Value = -3747

# row 122
$ws.Range("H122").Value = 962.75
$ws.Range("J122").Value = 1025.4
$ws.Range("L122").Value = 9228.6
$ws.Range("N122").Value = -14128.6

$ws = $wb.Worksheets.Item("GSM")
# row 10
$ws.Range("H10").Value = 9150
$ws.Range("I10").Value = 9150
$ws.Range("K10").Value = 9150
$ws.Range("M10").Value = -8981

# row 132
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
# row 10
$ws.Range("H10").Value = 1000000
$ws.Range("I10").Value = 1000000
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 1000000
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -999860
$ws.Range("N10").ClearContents()

# row 40
$ws.Range("H40").Value = 5000
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()

# row 46
$ws.Range("H46").Value = 58662.668
$ws.Range("I46").Value = 101992.8
$ws.Range("J46").Value = 4500
$ws.Range("K46").Value = 101992.8
$ws.Range("L46").Value = 4500
$ws.Range("M46").Value = -101804.8
$ws.Range("N46").Value = -4876

# row 61
$ws.Range("H61").Value = 3886.25
$ws.Range("I61").Value = 3182.5
$ws.Range("J61").Value = 5997.5
$ws.Range("K61").Value = 3182.5
$ws.Range("L61").Value = 5997.5
$ws.Range("M61").Value = -2980.5
$ws.Range("N61").Value = -6401.5

# row 111
$ws.Range("H111").Value = 40000
$ws.Range("J111").Value = 40000
$ws.Range("L111").Value = 40000
$ws.Range("N111").Value = -48180

# row 113
$ws.Range("H113").Value = 3886.25
$ws.Range("I113").Value = 3182.5
$ws.Range("J113").Value = 5997.5
$ws.Range("K113").Value = 3182.5
$ws.Range("L113").Value = 5997.5
$ws.Range("M113").Value = -1012.5
$ws.Range("N113").Value = -10337.5

# row 132
$ws.Range("H132").Value = 4746.2856
$ws.Range("I132").Value = 4746.2856
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 14238.8568
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -11708.8568
$ws.Range("N132").ClearContents()

# row 139
$ws.Range("H139").Value = 30000
$ws.Range("I139").Value = 30000
$ws.Range("K139").Value = 30000
$ws.Range("M139").Value = -24860

$ws = $wb.Worksheets.Item("WVR")
# row 37
$ws.Range("H37").Value = 20029
$ws.Range("J37").Value = 20029
$ws.Range("L37").Value = 20029
$ws.Range("N37").Value = -20435

# row 100
$ws.Range("H100").Value = 33334424
$ws.Range("I100").Value = 33334424
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 66668848
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -66668307
$ws.Range("N100").ClearContents()

# row 132
$ws.Range("H132").Value = 5699.391
$ws.Range("I132").Value = 5699.391
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 17098.173
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -14568.173
$ws.Range("N132").ClearContents()
